$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D1").Value = "sex"
